$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Posted At" holds date-look-alike strings ("2026-02-20") that must stay
# as plain text (matching the source inlineStr cells) instead of being
# auto-converted to Excel date serials. Force a text number format right
# before writing each one, then strip the format back off so no stray
# style is left behind on the cell.
function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range($cellRef).ClearFormats()
}

# Update existing row 2
$ws.Range("A2").Value = "Senior Software Engineer"
$ws.Range("B2").Value = "HCA Healthcare"
$ws.Range("C2").Value = "Nashville, TN, US USA"
$ws.Range("D2").Value = 11.1
$ws.Range("E2").Value = "RAG, Copilot, Docker, Kubernetes, AKS, Git, Python, SQL, R, Scala"
Set-TextValue "F2" "2026-02-20"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=7055f93be6329f9d"

# Update existing row 3
$ws.Range("A3").Value = "Analytics Developer 3 - HEDIS Quality"
$ws.Range("B3").Value = "Baylor Scott & White Health"
$ws.Range("C3").Value = "Remote, US USA"
$ws.Range("D3").Value = 11.1
$ws.Range("E3").Value = "Data Scientist, RAG, Cortex, Snowflake, Databricks, Power BI, Python, SQL, R, Scala"
Set-TextValue "F3" "2026-02-20"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=49c56062f945cbd6"

# New row 4
$ws.Range("A4").Value = "Software Engineer"
$ws.Range("B4").Value = "Forvia"
$ws.Range("C4").Value = "Auburn Hills, MI, US USA"
$ws.Range("D4").Value = 11.1
$ws.Range("E4").Value = "Generative AI, RAG, Docker, Kubernetes, CI/CD, Terraform, Python, R, Java, Scala"
Set-TextValue "F4" "2026-02-19"
$ws.Range("G4").Value = "https://www.indeed.com/viewjob?jk=711075e72263a26e"

# New row 5
$ws.Range("A5").Value = "Data Scientist"
$ws.Range("B5").Value = "Forvia"
$ws.Range("C5").Value = "Auburn Hills, MI, US USA"
$ws.Range("D5").Value = 11.1
$ws.Range("E5").Value = "Data Scientist, Machine Learning Engineer, Generative AI, RAG, TensorFlow, PyTorch, Python, R, Scala, Optimization"
Set-TextValue "F5" "2026-02-19"
$ws.Range("G5").Value = "https://www.indeed.com/viewjob?jk=0761c4262d5d60ff"
